$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new log entry as row 8
# Force column A to stay plain text so the date-like string isn't
# auto-converted into a date serial number, then restore the cell's
# style to the plain default style (same as the other "A" cells) so
# no extra style is left behind.
$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = "11.14.19"
$ws.Range("A8").Style = $ws.Range("A7").Style
$ws.Range("B8").Value = "Notebook refactoring, initial algorithm evaluatiojn section writing"

# Move the active selection, matching the saved workbook state
$ws.Range("B11").Select()
